# Updated cryptos list with GitHub Actions
# Applies latest scraped price/volume data to the crypto ranking sheet,
# including a text-formatted "Price" column (column D) and the
# percentage-styled "Volume(1h)" column (column E). Two rows (OKB/ONDO)
# also swap rank position as their underlying data order changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values must stay as plain text so values such as
# "1.00" or "0.999" keep their exact displayed digits instead of being
# coerced into numbers by Excel's automatic type detection.
$priceUpdates = @{
    "D2" = '63.560.95'
    "D3" = '3.484.07'
    "D4" = '1.00'
    "D5" = '581.20'
    "D6" = '147.89'
    "D7" = '3.483.60'
    "D10" = '7.72'
    "D13" = '4.078.42'
    "D14" = '29.73'
    "D16" = '3.473.81'
    "D18" = '63.517.02'
    "D21" = '9.37'
    "D22" = '391.29'
    "D23" = '0.564'
    "D24" = '75.04'
    "D26" = '3.612.19'
    "D30" = '0.999'
    "D35" = '23.65'
    "D36" = '7.19'
    "D37" = '32.49'
    "D38" = '5.33'
    "D40" = '170.28'
    "D41" = '3.519.75'
    "D44" = '1.23'
    "D45" = '42.46'
    "D48" = '2.628.84'
    "D50" = '23.20'
    "D51" = '6.80'
}

foreach ($cellRef in $priceUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$cellRef]
}

# Remaining updates: Coin name (B), Link (C) and Volume(1h) (E) columns.
$otherUpdates = @{
    "E2" = '  +3.01%  '
    "E3" = '  +2.12%  '
    "E4" = '  +0.07%  '
    "E5" = '  +0.72%  '
    "E6" = '  +2.77%  '
    "E7" = '  +2.10%  '
    "E8" = '  -0.16%  '
    "E9" = '  +1.09%  '
    "E10" = '  +1.14%  '
    "E11" = '  +2.24%  '
    "E12" = '  +5.32%  '
    "E13" = '  +2.12%  '
    "E14" = '  +6.26%  '
    "E15" = '  +2.74%  '
    "E16" = '  +1.64%  '
    "E17" = '  +1.69%  '
    "E18" = '  +2.88%  '
    "E19" = '  +3.57%  '
    "E20" = '  +5.41%  '
    "E21" = '  +2.17%  '
    "E22" = '  +0.82%  '
    "E23" = '  +2.50%  '
    "E24" = '  +1.04%  '
    "E25" = '  +0.09%  '
    "E26" = '  +1.75%  '
    "E27" = '  +1.35%  '
    "E28" = '  -3.09%  '
    "E29" = '  +3.36%  '
    "E30" = '  -0.01%  '
    "E31" = '  +3.13%  '
    "E32" = '  -0.52%  '
    "E33" = '  -0.03%  '
    "E34" = '  -2.37%  '
    "E35" = '  +0.93%  '
    "E36" = '  +3.52%  '
    "E37" = '  +14.68%  '
    "E38" = '  +2.97%  '
    "E39" = '  +7.82%  '
    "E40" = '  +1.19%  '
    "E41" = '  +2.16%  '
    "E42" = '  +1.37%  '
    "E43" = '  +1.90%  '
    "B44" = 'ONDO'
    "C44" = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
    "E44" = '  +5.43%  '
    "B45" = 'OKB'
    "C45" = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    "E45" = '  +0.02%  '
    "E46" = '  +3.56%  '
    "E47" = '  -0.17%  '
    "E48" = '  +5.58%  '
    "E49" = '  +10.96%  '
    "E50" = '  +1.70%  '
    "E51" = '  +2.56%  '
}

foreach ($cellRef in $otherUpdates.Keys) {
    $ws.Range($cellRef).Value = $otherUpdates[$cellRef]
}
